$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update relative influence values and reorder SSTmax/Salinity rows
$ws.Range("B2").Value = 38.98305797969099

$ws.Range("A3").Value = "SSTmax"
$ws.Range("B3").Value = 22.09628308700998

$ws.Range("A4").Value = "Salinity"
$ws.Range("B4").Value = 21.99824228171086

$ws.Range("B5").Value = 16.92241665158817
